$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autoavaliação")

# Identification
$ws.Range("B4").Value = "RCL10"

# Final grade
$ws.Range("B5").Value = 75

# Cliente TCP (left column B) and Cliente UDP (right column F) scores
$ws.Range("B9").Value = 1
$ws.Range("F9").Value = 0.25

$ws.Range("B10").Value = 0.5
$ws.Range("F10").Value = 0.5

$ws.Range("B11").Value = 0.5
$ws.Range("F11").Value = 0.5

$ws.Range("B12").Value = 1
$ws.Range("F12").Value = 0.75

$ws.Range("B13").Value = 1
$ws.Range("F13").Value = 1

$ws.Range("B14").Value = 0.25
$ws.Range("F14").Value = 0.25

$ws.Range("B15").Value = 0.25
$ws.Range("F15").Value = 0.25

$ws.Range("B17").Value = 1.5
$ws.Range("F17").Value = 1

$ws.Range("B21").Value = 2
$ws.Range("F21").Value = 0.75

$ws.Range("B22").Value = 0.5
$ws.Range("F22").Value = 0.5

$ws.Range("B23").Value = 1
$ws.Range("F23").Value = 0.75

$ws.Range("B24").Value = 0.5
$ws.Range("F24").Value = 0.5

$ws.Range("B26").Value = 2
$ws.Range("F26").Value = 0.75

# Script execution results: mark all as "Completo"
for ($row = 32; $row -le 55; $row++) {
    $ws.Range("B$row").Value = "Completo"
}

# Leave the selection where the author last left it
$ws.Range("G30").Select() | Out-Null
